$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add the new "caseid" column (I) ---
$ws.Range("I1").Value = "caseid"

# --- Row 2: now Carlos' record (was Ana's) ---
$ws.Range("A2").Value = "whatsapp:+573206198169"
$ws.Range("B2").Value = "Carlos"
$ws.Range("C2").Value = "Carlos Bohm"
$ws.Range("D2").Value = "diciembre"
$ws.Range("H2").Value = "working"
$ws.Range("I2").Value = "TestFelipe"

# --- Row 3: now Felipe's record (was Felipe's, but different phone/full name) ---
$ws.Range("A3").Value = "whatsapp:+17733220947"
$ws.Range("B3").Value = "Felipe"
$ws.Range("C3").Value = "Felipe Alamos"
$ws.Range("D3").Value = "enero"
$ws.Range("H3").Value = "looking for a job"
$ws.Range("I3").Value = "TestFelipe"

# --- Remove the fourth data row (Kiara) entirely - only two respondents remain ---
$ws.Rows.Item(4).Delete()

# --- Column sizing to fit the new, narrower data set ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(9).AutoFit()

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("E6").Select()
